# Apply cryptocurrency price/volume updates to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "317.58") must be
# forced to Text format first, otherwise Excel auto-converts the typed
# string into a numeric value, losing the original text cell type that
# the source data used (prices are stored as text, not numbers).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Write the updated price / volume text values.
$ws.Range('D2').Value = '42.745.83'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '2.532.85'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '317.58'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '97.48'
$ws.Range('E6').Value = '  +2.02%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').Value = '35.88'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('E13').Value = '  -2.44%  '
$ws.Range('D14').Value = '2.919.81'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').Value = '2.535.36'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '15.12'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').Value = '0.851'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '42.804.84'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('E19').Value = '  +4.44%  '
$ws.Range('D20').Value = '12.79'
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = '69.72'
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').Value = '251.65'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').Value = '26.48'
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').Value = '40.99'
$ws.Range('E29').Value = '  +5.75%  '
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').Value = '158.91'
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('D34').Value = '2.73'
$ws.Range('E34').Value = '  +4.53%  '
$ws.Range('D35').Value = '3.36'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').Value = '18.90'
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('D37').Value = '0.0788'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').Value = '2.33'
$ws.Range('E40').Value = '  +10.33%  '
$ws.Range('D41').Value = '22.31'
$ws.Range('E41').Value = '  -6.85%  '
$ws.Range('D42').Value = '3.83'
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  -2.99%  '
$ws.Range('D46').Value = '2.027.64'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('E47').Value = '  +3.27%  '
$ws.Range('D48').Value = '84.44'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '106.33'
$ws.Range('E49').Value = '  +4.94%  '
$ws.Range('D50').Value = '75.19'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').Value = '2.772.69'
$ws.Range('E51').Value = '  +0.63%  '
